$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.032238159739954
$ws.Cells.Item(2, 4).Value = 1.035411514214439
$ws.Cells.Item(2, 5).Value = 1.040941993116917
$ws.Cells.Item(2, 6).Value = 1.050074940190408
$ws.Cells.Item(2, 9).Value = 1.033533150488972
$ws.Cells.Item(2, 10).Value = 1.037369480076597
$ws.Cells.Item(2, 11).Value = 1.038208195315353
$ws.Cells.Item(2, 12).Value = 1.043722913166294
$ws.Cells.Item(2, 13).Value = 1.052830226853501

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033222425333232
$ws.Cells.Item(3, 4).Value = 1.036134189439881
$ws.Cells.Item(3, 5).Value = 1.041823143574659
$ws.Cells.Item(3, 6).Value = 1.051055720509156
$ws.Cells.Item(3, 9).Value = 1.033693886621651
$ws.Cells.Item(3, 10).Value = 1.037995779802559
$ws.Cells.Item(3, 11).Value = 1.03874049049689
$ws.Cells.Item(3, 12).Value = 1.044414404366428
$ws.Cells.Item(3, 13).Value = 1.05362294654327

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033859774066381
$ws.Cells.Item(4, 4).Value = 1.036601951805587
$ws.Cells.Item(4, 5).Value = 1.042394089855593
$ws.Cells.Item(4, 6).Value = 1.051691228812076
$ws.Cells.Item(4, 9).Value = 1.033796635727209
$ws.Cells.Item(4, 10).Value = 1.038400907677486
$ws.Cells.Item(4, 11).Value = 1.039084395181607
$ws.Cells.Item(4, 12).Value = 1.044861991377187
$ws.Cells.Item(4, 13).Value = 1.054136149083096

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034127825593162
$ws.Cells.Item(5, 4).Value = 1.036798632086645
$ws.Cells.Item(5, 5).Value = 1.042634301950497
$ws.Cells.Item(5, 6).Value = 1.051958605483189
$ws.Cells.Item(5, 9).Value = 1.033839529804051
$ws.Cells.Item(5, 10).Value = 1.038571191472869
$ws.Cells.Item(5, 11).Value = 1.039228845869806
$ws.Cells.Item(5, 12).Value = 1.045050190956526
$ws.Cells.Item(5, 13).Value = 1.054351960606277

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.034172839059974
$ws.Cells.Item(6, 4).Value = 1.036831657467779
$ws.Cells.Item(6, 5).Value = 1.042674645507492
$ws.Cells.Item(6, 6).Value = 1.052003511423274
$ws.Cells.Item(6, 9).Value = 1.033846714205096
$ws.Cells.Item(6, 10).Value = 1.038599780972136
$ws.Cells.Item(6, 11).Value = 1.039253092310567
$ws.Cells.Item(6, 12).Value = 1.045081792460367
$ws.Cells.Item(6, 13).Value = 1.054388199867703

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.033863355352226
$ws.Cells.Item(7, 4).Value = 1.036604579729156
$ws.Cells.Item(7, 5).Value = 1.042397298850889
$ws.Cells.Item(7, 6).Value = 1.051694800692567
$ws.Cells.Item(7, 9).Value = 1.033797210065851
$ws.Cells.Item(7, 10).Value = 1.038403183144562
$ws.Cells.Item(7, 11).Value = 1.039086325838355
$ws.Cells.Item(7, 12).Value = 1.044864505976764
$ws.Cells.Item(7, 13).Value = 1.05413903252816

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.032570700574357
$ws.Cells.Item(8, 4).Value = 1.035655715268952
$ws.Cells.Item(8, 5).Value = 1.041239619203689
$ws.Cells.Item(8, 6).Value = 1.05040621701599
$ws.Cells.Item(8, 9).Value = 1.033587732178368
$ws.Cells.Item(8, 10).Value = 1.037581167387215
$ws.Cells.Item(8, 11).Value = 1.038388195287138
$ws.Cells.Item(8, 12).Value = 1.043956575063674
$ws.Cells.Item(8, 13).Value = 1.053098075743804

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030296457052601
$ws.Cells.Item(9, 4).Value = 1.033984848373055
$ws.Cells.Item(9, 5).Value = 1.039205688612879
$ws.Cells.Item(9, 6).Value = 1.048142343295738
$ws.Cells.Item(9, 9).Value = 1.033208993213326
$ws.Cells.Item(9, 10).Value = 1.036131713349781
$ws.Cells.Item(9, 11).Value = 1.037154013877927
$ws.Cells.Item(9, 12).Value = 1.042357844361874
$ws.Cells.Item(9, 13).Value = 1.051265813390655

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028782743725124
$ws.Cells.Item(10, 4).Value = 1.032871788747539
$ws.Cells.Item(10, 5).Value = 1.037853870803772
$ws.Cells.Item(10, 6).Value = 1.046637718527403
$ws.Cells.Item(10, 9).Value = 1.032950067206241
$ws.Cells.Item(10, 10).Value = 1.035164817161202
$ws.Cells.Item(10, 11).Value = 1.036328601411075
$ws.Cells.Item(10, 12).Value = 1.041292858766433
$ws.Cells.Item(10, 13).Value = 1.050045735426713

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.028127878424783
$ws.Cells.Item(11, 4).Value = 1.032390039436615
$ws.Cells.Item(11, 5).Value = 1.037269513991818
$ws.Cells.Item(11, 6).Value = 1.045987310412934
$ws.Cells.Item(11, 9).Value = 1.032836429526944
$ws.Cells.Item(11, 10).Value = 1.034746009930003
$ws.Cells.Item(11, 11).Value = 1.035970577510792
$ws.Cells.Item(11, 12).Value = 1.040831917618992
$ws.Cells.Item(11, 13).Value = 1.049517779979842

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.02788472042617
$ws.Cells.Item(12, 4).Value = 1.032211129514219
$ws.Cells.Item(12, 5).Value = 1.037052607586331
$ws.Cells.Item(12, 6).Value = 1.045745886895775
$ws.Cells.Item(12, 9).Value = 1.032793991374344
$ws.Cells.Item(12, 10).Value = 1.034590426810188
$ws.Cells.Item(12, 11).Value = 1.035837499997514
$ws.Cells.Item(12, 12).Value = 1.040660735175556
$ws.Cells.Item(12, 13).Value = 1.049321726825228

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02793687464431
$ws.Cells.Item(13, 4).Value = 1.032249504776604
$ws.Cells.Item(13, 5).Value = 1.037099127952986
$ws.Cells.Item(13, 6).Value = 1.045797665484201
$ws.Cells.Item(13, 9).Value = 1.032803104814947
$ws.Cells.Item(13, 10).Value = 1.03462380078581
$ws.Cells.Item(13, 11).Value = 1.035866049702102
$ws.Cells.Item(13, 12).Value = 1.040697452945705
$ws.Cells.Item(13, 13).Value = 1.049363778478226

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.02810777710736
$ws.Cells.Item(14, 4).Value = 1.032375250001724
$ws.Cells.Item(14, 5).Value = 1.037251581379723
$ws.Cells.Item(14, 6).Value = 1.045967350856238
$ws.Cells.Item(14, 9).Value = 1.032832926225702
$ws.Cells.Item(14, 10).Value = 1.034733149773351
$ws.Cells.Item(14, 11).Value = 1.035959579143427
$ws.Cells.Item(14, 12).Value = 1.040817766974658
$ws.Cells.Item(14, 13).Value = 1.049501573078016

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.028213087451969
$ws.Cells.Item(15, 4).Value = 1.032452730224379
$ws.Cells.Item(15, 5).Value = 1.037345532841941
$ws.Cells.Item(15, 6).Value = 1.046071921778504
$ws.Cells.Item(15, 9).Value = 1.03285126997475
$ws.Cells.Item(15, 10).Value = 1.034800520737539
$ws.Cells.Item(15, 11).Value = 1.036017193621335
$ws.Cells.Item(15, 12).Value = 1.040891900635253
$ws.Cells.Item(15, 13).Value = 1.049586479931501

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028826217304063
$ws.Cells.Item(16, 4).Value = 1.032903765453369
$ws.Cells.Item(16, 5).Value = 1.037892673577626
$ws.Cells.Item(16, 6).Value = 1.046680907333017
$ws.Cells.Item(16, 9).Value = 1.032957576957702
$ws.Cells.Item(16, 10).Value = 1.035192609244745
$ws.Cells.Item(16, 11).Value = 1.036352349398967
$ws.Cells.Item(16, 12).Value = 1.041323454298722
$ws.Cells.Item(16, 13).Value = 1.050080781471028

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029210973857526
$ws.Cells.Item(17, 4).Value = 1.033186745800674
$ws.Cells.Item(17, 5).Value = 1.038236146134155
$ws.Cells.Item(17, 6).Value = 1.047063204322041
$ws.Cells.Item(17, 9).Value = 1.033023853623361
$ws.Cells.Item(17, 10).Value = 1.035438520470568
$ws.Cells.Item(17, 11).Value = 1.036562419770016
$ws.Cells.Item(17, 12).Value = 1.041594211917702
$ws.Cells.Item(17, 13).Value = 1.050390937181769

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029435452013017
$ws.Cells.Item(18, 4).Value = 1.033351823903438
$ws.Cells.Item(18, 5).Value = 1.038436583235917
$ws.Cells.Item(18, 6).Value = 1.047286298166507
$ws.Cells.Item(18, 9).Value = 1.033062364895124
$ws.Cells.Item(18, 10).Value = 1.035581943268992
$ws.Cells.Item(18, 11).Value = 1.036684890879716
$ws.Cells.Item(18, 12).Value = 1.041752159942221
$ws.Cells.Item(18, 13).Value = 1.05057187911435

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029512002771529
$ws.Cells.Item(19, 4).Value = 1.033408114701976
$ws.Cells.Item(19, 5).Value = 1.038504943250503
$ws.Cells.Item(19, 6).Value = 1.047362385429373
$ws.Cells.Item(19, 9).Value = 1.033075471332279
$ws.Cells.Item(19, 10).Value = 1.035630844497951
$ws.Cells.Item(19, 11).Value = 1.036726640263953
$ws.Cells.Item(19, 12).Value = 1.041806019449312
$ws.Cells.Item(19, 13).Value = 1.050633581265374

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029169687306525
$ws.Cells.Item(20, 4).Value = 1.033156382562867
$ws.Cells.Item(20, 5).Value = 1.038199284906722
$ws.Cells.Item(20, 6).Value = 1.047022176451553
$ws.Cells.Item(20, 9).Value = 1.033016757943835
$ws.Cells.Item(20, 10).Value = 1.035412137868376
$ws.Cells.Item(20, 11).Value = 1.036539887343236
$ws.Cells.Item(20, 12).Value = 1.041565160145006
$ws.Cells.Item(20, 13).Value = 1.050357656952307

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.028057448153838
$ws.Cells.Item(21, 4).Value = 1.032338220239989
$ws.Cells.Item(21, 5).Value = 1.037206683503037
$ws.Cells.Item(21, 6).Value = 1.045917378124153
$ws.Cells.Item(21, 9).Value = 1.032824150855386
$ws.Cells.Item(21, 10).Value = 1.034700949753033
$ws.Cells.Item(21, 11).Value = 1.035932039572781
$ws.Cells.Item(21, 12).Value = 1.040782336611906
$ws.Cells.Item(21, 13).Value = 1.049460994531816

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.027358649462838
$ws.Cells.Item(22, 4).Value = 1.031824001939178
$ws.Cells.Item(22, 5).Value = 1.036583462465288
$ws.Cells.Item(22, 6).Value = 1.04522371441645
$ws.Cells.Item(22, 9).Value = 1.03270173188306
$ws.Cells.Item(22, 10).Value = 1.034253685319358
$ws.Cells.Item(22, 11).Value = 1.035549332623949
$ws.Cells.Item(22, 12).Value = 1.040290327777487
$ws.Cells.Item(22, 13).Value = 1.04889753443449

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.02772904730401
$ws.Cells.Item(23, 4).Value = 1.032096580097479
$ws.Cells.Item(23, 5).Value = 1.036913761188507
$ws.Cells.Item(23, 6).Value = 1.045591346657798
$ws.Cells.Item(23, 9).Value = 1.032766753393796
$ws.Cells.Item(23, 10).Value = 1.034490799034381
$ws.Cells.Item(23, 11).Value = 1.035752262661277
$ws.Cells.Item(23, 12).Value = 1.040551133293232
$ws.Cells.Item(23, 13).Value = 1.049196205854447

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029188342748204
$ws.Cells.Item(24, 4).Value = 1.033170102339976
$ws.Cells.Item(24, 5).Value = 1.038215940614136
$ws.Cells.Item(24, 6).Value = 1.047040714851277
$ws.Cells.Item(24, 9).Value = 1.033019964629561
$ws.Cells.Item(24, 10).Value = 1.035424059071078
$ws.Cells.Item(24, 11).Value = 1.036550068960551
$ws.Cells.Item(24, 12).Value = 1.041578287329801
$ws.Cells.Item(24, 13).Value = 1.050372694752477

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030883974986473
$ws.Cells.Item(25, 4).Value = 1.034416662974905
$ws.Cells.Item(25, 5).Value = 1.039730784495487
$ws.Cells.Item(25, 6).Value = 1.048726798899577
$ws.Cells.Item(25, 9).Value = 1.033308042388879
$ws.Cells.Item(25, 10).Value = 1.036506540242277
$ws.Cells.Item(25, 11).Value = 1.037473546049255
$ws.Cells.Item(25, 12).Value = 1.042771011466525
$ws.Cells.Item(25, 13).Value = 1.051739249508673
